$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(12)
$shape.Top = 7.66906
